$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column H (TEST), shifting H:R -> I:S.
# Restrict the insert to just the used rows (1:11) so column-level
# formatting metadata for untouched columns is left alone.
$ws.Range("H1:H11").Insert(-4161)

# New header for the inserted column, matching the look of its neighbor (old TEST header, now I1).
$ws.Range("I1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "#TUS"

# Give the new data cells (H2:H10) the same formatting as their neighbor column (I2:I10).
$ws.Range("I2:I10").Copy()
$ws.Range("H2:H10").PasteSpecial(-4122)

# Populate the new #TUS column with the value 5 for every data row.
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 8).Value = 5
}

$excel.CutCopyMode = $false
